$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the first empty row right after the existing data (row 6 in this case)
$newRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1

$ws.Cells.Item($newRow, 1).Value = "JD_005"
$ws.Cells.Item($newRow, 2).Value = "Dummy Engineer"
$ws.Cells.Item($newRow, 3).Value = "Dummy Data"
$ws.Cells.Item($newRow, 4).Value = 1
$ws.Cells.Item($newRow, 5).Value = 5
